$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update custom_fields_json (column M) for rows 39-47 ---
$ws.Range("M39").Value = "{`"custom_1`": `"ParcelHub`", `"custom_2`": `"TRL-MER`"}"
$ws.Range("M40").Value = "{`"custom_1`": `"ParcelHub`", `"custom_2`": `"TRL-CLE`"}"
$ws.Range("M42").Value = "{`"custom_1`": `"ParcelHub`", `"custom_2`": `"TRL-RES`"}"
$ws.Range("M43").Value = "{`"custom_1`": `"ParcelHub`", `"custom_2`": `"TRL-RES`"}"
$ws.Range("M44").Value = "{`"custom_1`": `"ParcelHub`", `"custom_2`": `"TRL-MER`"}"
$ws.Range("M45").Value = "{`"custom_1`": `"ParcelHub`", `"custom_2`": `"TRL-RES`"}"
$ws.Range("M47").Value = "{`"custom_1`": `"ParcelHub`", `"custom_2`": `"TRL-MER`"}"

# --- Reassign rows 145-201 (permutation of existing shipment rows) ---
# row 145 <- old row 148
$ws.Range("A145").Value = "JJD149990121180193812"
$ws.Range("B145").Value = "dhlparcel-nl"
$ws.Range("C145").Value = "DHL Parcel NL"
$ws.Range("I145").Value = "2026-02-04T14:34:00+01:00"
$ws.Range("J145").Value = ""
$ws.Range("K145").Value = "2026-02-09T01:39:10+00:00"
$ws.Range("L145").Value = "JJD149990121180193812"

# row 146 <- old row 150
$ws.Range("A146").Value = "JJD149990121180194480"
$ws.Range("B146").Value = "dhlparcel-nl"
$ws.Range("C146").Value = "DHL Parcel NL"
$ws.Range("I146").Value = "2026-02-05T11:49:00+01:00"
$ws.Range("J146").Value = ""
$ws.Range("K146").Value = "2026-02-09T01:39:10+00:00"
$ws.Range("L146").Value = "JJD149990121180194480"

# row 147 <- old row 152
$ws.Range("A147").Value = "JJD149990121180193693"
$ws.Range("B147").Value = "dhlparcel-nl"
$ws.Range("C147").Value = "DHL Parcel NL"
$ws.Range("I147").Value = "2026-02-04T14:34:00+01:00"
$ws.Range("J147").Value = ""
$ws.Range("K147").Value = "2026-02-09T01:39:08+00:00"
$ws.Range("L147").Value = "JJD149990121180193693"

# row 148 <- old row 166
$ws.Range("A148").NumberFormat = "@"
$ws.Range("A148").Value = "09447272602648"
$ws.Range("B148").Value = "dpd"
$ws.Range("C148").Value = "DPD"
$ws.Range("I148").Value = "2026-01-16T14:43:03+01:00"
$ws.Range("J148").Value = "Aschaffenburg, DE, Germany"
$ws.Range("K148").Value = "2026-02-09T01:39:11+00:00"
$ws.Range("L148").NumberFormat = "@"
$ws.Range("L148").Value = "09447272602648"

# row 150 <- old row 156
$ws.Range("A150").Value = "JJD149990121180193851"
$ws.Range("B150").Value = "dhlparcel-nl"
$ws.Range("C150").Value = "DHL Parcel NL"
$ws.Range("I150").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J150").Value = ""
$ws.Range("K150").Value = "2026-02-09T16:09:20+00:00"
$ws.Range("L150").Value = "JJD149990121180193851"

# row 151 <- old row 157
$ws.Range("A151").NumberFormat = "@"
$ws.Range("A151").Value = "01606828664171"
$ws.Range("B151").Value = "dpd-de"
$ws.Range("C151").Value = "DPD Germany"
$ws.Range("I151").Value = "2026-02-09T08:20:00+01:00"
$ws.Range("J151").Value = "Oirschot (NL), Netherlands"
$ws.Range("K151").Value = "2026-02-09T13:39:20+00:00"
$ws.Range("L151").NumberFormat = "@"
$ws.Range("L151").Value = "01606828664171"

# row 152 <- old row 155
$ws.Range("A152").Value = "JJD149990121180189785"
$ws.Range("B152").Value = "dhlparcel-nl"
$ws.Range("C152").Value = "DHL Parcel NL"
$ws.Range("I152").Value = "2026-01-30T16:38:00+01:00"
$ws.Range("J152").Value = ""
$ws.Range("K152").Value = "2026-02-09T01:39:11+00:00"
$ws.Range("L152").Value = "JJD149990121180189785"

# row 153 <- old row 154
$ws.Range("A153").Value = "JJD149990121180194493"
$ws.Range("B153").Value = "dhlparcel-nl"
$ws.Range("C153").Value = "DHL Parcel NL"
$ws.Range("I153").Value = "2026-02-03T14:59:00+01:00"
$ws.Range("J153").Value = ""
$ws.Range("K153").Value = "2026-02-09T01:39:08+00:00"
$ws.Range("L153").Value = "JJD149990121180194493"

# row 154 <- old row 159
$ws.Range("A154").Value = "JJD149990121180193486"
$ws.Range("B154").Value = "dhlparcel-nl"
$ws.Range("C154").Value = "DHL Parcel NL"
$ws.Range("I154").Value = "2026-02-03T14:59:00+01:00"
$ws.Range("J154").Value = ""
$ws.Range("K154").Value = "2026-02-09T01:39:11+00:00"
$ws.Range("L154").Value = "JJD149990121180193486"

# row 155 <- old row 162
$ws.Range("A155").Value = "JJD149990121180193821"
$ws.Range("B155").Value = "dhlparcel-nl"
$ws.Range("C155").Value = "DHL Parcel NL"
$ws.Range("I155").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J155").Value = ""
$ws.Range("K155").Value = "2026-02-09T16:09:32+00:00"
$ws.Range("L155").Value = "JJD149990121180193821"

# row 156 <- old row 146
$ws.Range("A156").Value = "JJD149990121180189756"
$ws.Range("B156").Value = "dhlparcel-nl"
$ws.Range("C156").Value = "DHL Parcel NL"
$ws.Range("I156").Value = "2026-02-03T14:59:00+01:00"
$ws.Range("J156").Value = ""
$ws.Range("K156").Value = "2026-02-09T01:39:10+00:00"
$ws.Range("L156").Value = "JJD149990121180189756"

# row 157 <- old row 153
$ws.Range("A157").Value = "JJD149990121180192273"
$ws.Range("B157").Value = "dhlparcel-nl"
$ws.Range("C157").Value = "DHL Parcel NL"
$ws.Range("I157").Value = "2026-02-02T14:31:00+01:00"
$ws.Range("J157").Value = ""
$ws.Range("K157").Value = "2026-02-09T01:39:10+00:00"
$ws.Range("L157").Value = "JJD149990121180192273"

# row 158 <- old row 161
$ws.Range("A158").NumberFormat = "@"
$ws.Range("A158").Value = "09447272788603"
$ws.Range("B158").Value = "dpd"
$ws.Range("C158").Value = "DPD"
$ws.Range("I158").Value = "2026-02-09T08:20:00+01:00"
$ws.Range("J158").Value = "Oirschot (NL), Netherlands"
$ws.Range("K158").Value = "2026-02-09T14:39:14+00:00"
$ws.Range("L158").NumberFormat = "@"
$ws.Range("L158").Value = "09447272788603"

# row 159 <- old row 145
$ws.Range("A159").Value = "JJD149990121180193245"
$ws.Range("B159").Value = "dhlparcel-nl"
$ws.Range("C159").Value = "DHL Parcel NL"
$ws.Range("I159").Value = "2026-02-02T14:31:00+01:00"
$ws.Range("J159").Value = ""
$ws.Range("K159").Value = "2026-02-09T01:39:09+00:00"
$ws.Range("L159").Value = "JJD149990121180193245"

# row 160 <- old row 167
$ws.Range("A160").NumberFormat = "@"
$ws.Range("A160").Value = "166050916804958"
$ws.Range("B160").Value = "brt-it-parcelid"
$ws.Range("C160").Value = "BRT Bartolini(Parcel ID)"
$ws.Range("I160").Value = "2026-02-05T12:00:00+01:00"
$ws.Range("J160").Value = "MILANO SEDRIANO (050)"
$ws.Range("K160").Value = "2026-02-09T02:39:20+00:00"
$ws.Range("L160").NumberFormat = "@"
$ws.Range("L160").Value = "166050916804958"

# row 161 <- old row 164
$ws.Range("A161").Value = "JJD149990121180194786"
$ws.Range("B161").Value = "dhlparcel-nl"
$ws.Range("C161").Value = "DHL Parcel NL"
$ws.Range("I161").Value = "2026-02-03T14:59:00+01:00"
$ws.Range("J161").Value = ""
$ws.Range("K161").Value = "2026-02-09T01:39:10+00:00"
$ws.Range("L161").Value = "JJD149990121180194786"

# row 162 <- old row 151
$ws.Range("A162").Value = "JJD149990121180194455"
$ws.Range("B162").Value = "dhlparcel-nl"
$ws.Range("C162").Value = "DHL Parcel NL"
$ws.Range("I162").Value = "2026-02-05T11:49:00+01:00"
$ws.Range("J162").Value = ""
$ws.Range("K162").Value = "2026-02-09T01:39:09+00:00"
$ws.Range("L162").Value = "JJD149990121180194455"

# row 163 <- old row 158
$ws.Range("A163").Value = "CY425450001DE"
$ws.Range("B163").Value = "dhl-germany"
$ws.Range("C163").Value = "Deutsche Post DHL"
$ws.Range("I163").Value = "2026-01-22T14:45:00+01:00"
$ws.Range("J163").Value = "Netherlands"
$ws.Range("K163").Value = "2026-02-09T01:39:08+00:00"
$ws.Range("L163").Value = "CY425450001DE"

# row 164 <- old row 163
$ws.Range("A164").NumberFormat = "@"
$ws.Range("A164").Value = "09447272763065"
$ws.Range("B164").Value = "dpd"
$ws.Range("C164").Value = "DPD"
$ws.Range("I164").Value = "2026-01-28T08:53:02+01:00"
$ws.Range("J164").Value = "Wuppertal, DE, Germany"
$ws.Range("K164").Value = "2026-02-09T01:39:11+00:00"
$ws.Range("L164").NumberFormat = "@"
$ws.Range("L164").Value = "09447272763065"

# row 165 <- old row 147
$ws.Range("A165").NumberFormat = "@"
$ws.Range("A165").Value = "08458093138698"
$ws.Range("B165").Value = "brt-it"
$ws.Range("C165").Value = "BRT Bartolini"
$ws.Range("I165").Value = "2026-02-05T12:00:00+01:00"
$ws.Range("J165").Value = "MILANO SEDRIANO (050)"
$ws.Range("K165").Value = "2026-02-09T04:39:10+00:00"
$ws.Range("L165").NumberFormat = "@"
$ws.Range("L165").Value = "08458093138698"

# row 166 <- old row 160
$ws.Range("A166").NumberFormat = "@"
$ws.Range("A166").Value = "09447272761687"
$ws.Range("B166").Value = "dpd"
$ws.Range("C166").Value = "DPD"
$ws.Range("I166").Value = "2026-02-02T15:56:13+01:00"
$ws.Range("J166").Value = "Nagold, DE, Germany"
$ws.Range("K166").Value = "2026-02-09T01:39:11+00:00"
$ws.Range("L166").NumberFormat = "@"
$ws.Range("L166").Value = "09447272761687"

# row 167 <- old row 165
$ws.Range("A167").Value = "JJD149990121180200629"
$ws.Range("B167").Value = "dhlparcel-nl"
$ws.Range("C167").Value = "DHL Parcel NL"
$ws.Range("I167").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J167").Value = ""
$ws.Range("K167").Value = "2026-02-09T16:09:19+00:00"
$ws.Range("L167").Value = "JJD149990121180200629"

# row 168 <- old row 184
$ws.Range("A168").Value = "JJD149990121180195647"
$ws.Range("B168").Value = "dhlparcel-nl"
$ws.Range("C168").Value = "DHL Parcel NL"
$ws.Range("I168").Value = "2026-02-06T12:47:00+01:00"
$ws.Range("J168").Value = ""
$ws.Range("K168").Value = "2026-02-09T13:39:09+00:00"
$ws.Range("L168").Value = "JJD149990121180195647"

# row 169 <- old row 171
$ws.Range("A169").Value = "JJD149990121180196222"
$ws.Range("B169").Value = "dhlparcel-nl"
$ws.Range("C169").Value = "DHL Parcel NL"
$ws.Range("I169").Value = "2026-02-05T11:49:00+01:00"
$ws.Range("J169").Value = ""
$ws.Range("K169").Value = "2026-02-09T01:39:08+00:00"
$ws.Range("L169").Value = "JJD149990121180196222"

# row 170 <- old row 182
$ws.Range("A170").Value = "JJD149990121180195426"
$ws.Range("B170").Value = "dhlparcel-nl"
$ws.Range("C170").Value = "DHL Parcel NL"
$ws.Range("I170").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J170").Value = ""
$ws.Range("K170").Value = "2026-02-09T16:09:18+00:00"
$ws.Range("L170").Value = "JJD149990121180195426"

# row 171 <- old row 189
$ws.Range("A171").Value = "JJD149990121180195416"
$ws.Range("B171").Value = "dhlparcel-nl"
$ws.Range("C171").Value = "DHL Parcel NL"
$ws.Range("I171").Value = "2026-02-04T14:34:00+01:00"
$ws.Range("J171").Value = ""
$ws.Range("K171").Value = "2026-02-09T01:39:09+00:00"
$ws.Range("L171").Value = "JJD149990121180195416"

# row 172 <- old row 192
$ws.Range("A172").Value = "JJD149990121180194932"
$ws.Range("B172").Value = "dhlparcel-nl"
$ws.Range("C172").Value = "DHL Parcel NL"
$ws.Range("I172").Value = "2026-02-06T12:47:00+01:00"
$ws.Range("J172").Value = ""
$ws.Range("K172").Value = "2026-02-09T13:39:10+00:00"
$ws.Range("L172").Value = "JJD149990121180194932"

# row 173 <- old row 193
$ws.Range("A173").Value = "JJD149990121180196232"
$ws.Range("B173").Value = "dhlparcel-nl"
$ws.Range("C173").Value = "DHL Parcel NL"
$ws.Range("I173").Value = "2026-02-04T14:34:00+01:00"
$ws.Range("J173").Value = ""
$ws.Range("K173").Value = "2026-02-09T01:39:07+00:00"
$ws.Range("L173").Value = "JJD149990121180196232"

# row 174 <- old row 173
$ws.Range("A174").Value = "JJD149990121180197789"
$ws.Range("B174").Value = "dhlparcel-nl"
$ws.Range("C174").Value = "DHL Parcel NL"
$ws.Range("I174").Value = "2026-02-05T11:49:00+01:00"
$ws.Range("J174").Value = ""
$ws.Range("K174").Value = "2026-02-09T01:39:06+00:00"
$ws.Range("L174").Value = "JJD149990121180197789"

# row 175 <- old row 176
$ws.Range("A175").Value = "JJD149990121180197333"
$ws.Range("B175").Value = "dhlparcel-nl"
$ws.Range("C175").Value = "DHL Parcel NL"
$ws.Range("I175").Value = "2026-02-05T11:49:00+01:00"
$ws.Range("J175").Value = ""
$ws.Range("K175").Value = "2026-02-09T01:39:08+00:00"
$ws.Range("L175").Value = "JJD149990121180197333"

# row 176 <- old row 190
$ws.Range("A176").Value = "JJD149990121180195903"
$ws.Range("B176").Value = "dhlparcel-nl"
$ws.Range("C176").Value = "DHL Parcel NL"
$ws.Range("I176").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J176").Value = ""
$ws.Range("K176").Value = "2026-02-09T16:09:18+00:00"
$ws.Range("L176").Value = "JJD149990121180195903"

# row 177 <- old row 174
$ws.Range("A177").Value = "JJD149990121180196585"
$ws.Range("B177").Value = "dhlparcel-nl"
$ws.Range("C177").Value = "DHL Parcel NL"
$ws.Range("I177").Value = "2026-02-05T11:49:00+01:00"
$ws.Range("J177").Value = ""
$ws.Range("K177").Value = "2026-02-09T01:39:07+00:00"
$ws.Range("L177").Value = "JJD149990121180196585"

# row 178 <- old row 179
$ws.Range("A178").Value = "JJD149990121180194902"
$ws.Range("B178").Value = "dhlparcel-nl"
$ws.Range("C178").Value = "DHL Parcel NL"
$ws.Range("I178").Value = "2026-02-04T14:34:00+01:00"
$ws.Range("J178").Value = ""
$ws.Range("K178").Value = "2026-02-09T01:39:09+00:00"
$ws.Range("L178").Value = "JJD149990121180194902"

# row 179 <- old row 168
$ws.Range("A179").Value = "JJD149990121180196112"
$ws.Range("B179").Value = "dhlparcel-nl"
$ws.Range("C179").Value = "DHL Parcel NL"
$ws.Range("I179").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J179").Value = ""
$ws.Range("K179").Value = "2026-02-09T16:09:21+00:00"
$ws.Range("L179").Value = "JJD149990121180196112"

# row 180 <- old row 169
$ws.Range("A180").Value = "JJD149990121180194983"
$ws.Range("B180").Value = "dhlparcel-nl"
$ws.Range("C180").Value = "DHL Parcel NL"
$ws.Range("I180").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J180").Value = ""
$ws.Range("K180").Value = "2026-02-09T16:09:20+00:00"
$ws.Range("L180").Value = "JJD149990121180194983"

# row 181 <- old row 180
$ws.Range("A181").Value = "JJD149990121180195907"
$ws.Range("B181").Value = "dhlparcel-nl"
$ws.Range("C181").Value = "DHL Parcel NL"
$ws.Range("I181").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J181").Value = ""
$ws.Range("K181").Value = "2026-02-09T16:09:32+00:00"
$ws.Range("L181").Value = "JJD149990121180195907"

# row 182 <- old row 178
$ws.Range("A182").Value = "JJD149990121180195374"
$ws.Range("B182").Value = "dhlparcel-nl"
$ws.Range("C182").Value = "DHL Parcel NL"
$ws.Range("I182").Value = "2026-02-04T14:34:00+01:00"
$ws.Range("J182").Value = ""
$ws.Range("K182").Value = "2026-02-09T01:39:08+00:00"
$ws.Range("L182").Value = "JJD149990121180195374"

# row 183 <- old row 181
$ws.Range("A183").Value = "JJD149990121180195611"
$ws.Range("B183").Value = "dhlparcel-nl"
$ws.Range("C183").Value = "DHL Parcel NL"
$ws.Range("I183").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J183").Value = ""
$ws.Range("K183").Value = "2026-02-09T16:09:21+00:00"
$ws.Range("L183").Value = "JJD149990121180195611"

# row 184 <- old row 187
$ws.Range("A184").Value = "JJD149990121180197142"
$ws.Range("B184").Value = "dhlparcel-nl"
$ws.Range("C184").Value = "DHL Parcel NL"
$ws.Range("I184").Value = "2026-02-05T11:49:00+01:00"
$ws.Range("J184").Value = ""
$ws.Range("K184").Value = "2026-02-09T01:39:09+00:00"
$ws.Range("L184").Value = "JJD149990121180197142"

# row 185 <- old row 177
$ws.Range("A185").Value = "JJD149990121180195413"
$ws.Range("B185").Value = "dhlparcel-nl"
$ws.Range("C185").Value = "DHL Parcel NL"
$ws.Range("I185").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J185").Value = ""
$ws.Range("K185").Value = "2026-02-09T16:09:18+00:00"
$ws.Range("L185").Value = "JJD149990121180195413"

# row 186 <- old row 183
$ws.Range("A186").Value = "JJD149990121180196685"
$ws.Range("B186").Value = "dhlparcel-nl"
$ws.Range("C186").Value = "DHL Parcel NL"
$ws.Range("I186").Value = "2026-02-06T12:47:00+01:00"
$ws.Range("J186").Value = ""
$ws.Range("K186").Value = "2026-02-09T13:39:08+00:00"
$ws.Range("L186").Value = "JJD149990121180196685"

# row 187 <- old row 175
$ws.Range("A187").Value = "JJD149990121180195036"
$ws.Range("B187").Value = "dhlparcel-nl"
$ws.Range("C187").Value = "DHL Parcel NL"
$ws.Range("I187").Value = "2026-02-04T14:34:00+01:00"
$ws.Range("J187").Value = ""
$ws.Range("K187").Value = "2026-02-09T01:39:09+00:00"
$ws.Range("L187").Value = "JJD149990121180195036"

# row 188 <- old row 185
$ws.Range("A188").Value = "JJD149990121180195981"
$ws.Range("B188").Value = "dhlparcel-nl"
$ws.Range("C188").Value = "DHL Parcel NL"
$ws.Range("I188").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J188").Value = ""
$ws.Range("K188").Value = "2026-02-09T16:09:21+00:00"
$ws.Range("L188").Value = "JJD149990121180195981"

# row 189 <- old row 172
$ws.Range("A189").Value = "JJD149990121180195961"
$ws.Range("B189").Value = "dhlparcel-nl"
$ws.Range("C189").Value = "DHL Parcel NL"
$ws.Range("I189").Value = "2026-02-04T14:34:00+01:00"
$ws.Range("J189").Value = ""
$ws.Range("K189").Value = "2026-02-09T01:39:08+00:00"
$ws.Range("L189").Value = "JJD149990121180195961"

# row 190 <- old row 188
$ws.Range("A190").Value = "JJD149990121180194833"
$ws.Range("B190").Value = "dhlparcel-nl"
$ws.Range("C190").Value = "DHL Parcel NL"
$ws.Range("I190").Value = "2026-02-06T12:47:00+01:00"
$ws.Range("J190").Value = ""
$ws.Range("K190").Value = "2026-02-09T13:39:12+00:00"
$ws.Range("L190").Value = "JJD149990121180194833"

# row 191 <- old row 186
$ws.Range("A191").Value = "JJD149990121180195737"
$ws.Range("B191").Value = "dhlparcel-nl"
$ws.Range("C191").Value = "DHL Parcel NL"
$ws.Range("I191").Value = "2026-02-05T11:49:00+01:00"
$ws.Range("J191").Value = ""
$ws.Range("K191").Value = "2026-02-09T01:39:09+00:00"
$ws.Range("L191").Value = "JJD149990121180195737"

# row 192 <- old row 191
$ws.Range("A192").Value = "JJD149990121180196977"
$ws.Range("B192").Value = "dhlparcel-nl"
$ws.Range("C192").Value = "DHL Parcel NL"
$ws.Range("I192").Value = "2026-02-04T14:34:00+01:00"
$ws.Range("J192").Value = ""
$ws.Range("K192").Value = "2026-02-09T01:39:08+00:00"
$ws.Range("L192").Value = "JJD149990121180196977"

# row 193 <- old row 170
$ws.Range("A193").Value = "JJD149990121180197491"
$ws.Range("B193").Value = "dhlparcel-nl"
$ws.Range("C193").Value = "DHL Parcel NL"
$ws.Range("I193").Value = "2026-02-06T12:47:00+01:00"
$ws.Range("J193").Value = ""
$ws.Range("K193").Value = "2026-02-09T13:39:09+00:00"
$ws.Range("L193").Value = "JJD149990121180197491"

# row 194 <- old row 199
$ws.Range("A194").Value = "JJD149990121180198565"
$ws.Range("B194").Value = "dhlparcel-nl"
$ws.Range("C194").Value = "DHL Parcel NL"
$ws.Range("I194").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J194").Value = ""
$ws.Range("K194").Value = "2026-02-09T16:09:15+00:00"
$ws.Range("L194").Value = "JJD149990121180198565"

# row 195 <- old row 197
$ws.Range("A195").Value = "JJD149990121180197838"
$ws.Range("B195").Value = "dhlparcel-nl"
$ws.Range("C195").Value = "DHL Parcel NL"
$ws.Range("I195").Value = "2026-02-05T11:49:00+01:00"
$ws.Range("J195").Value = ""
$ws.Range("K195").Value = "2026-02-09T01:39:07+00:00"
$ws.Range("L195").Value = "JJD149990121180197838"

# row 196 <- old row 198
$ws.Range("A196").Value = "JJD149990121180198097"
$ws.Range("B196").Value = "dhlparcel-nl"
$ws.Range("C196").Value = "DHL Parcel NL"
$ws.Range("I196").Value = "2026-02-06T12:47:00+01:00"
$ws.Range("J196").Value = ""
$ws.Range("K196").Value = "2026-02-09T13:39:08+00:00"
$ws.Range("L196").Value = "JJD149990121180198097"

# row 197 <- old row 200
$ws.Range("A197").Value = "JJD149990121180199217"
$ws.Range("B197").Value = "dhlparcel-nl"
$ws.Range("C197").Value = "DHL Parcel NL"
$ws.Range("I197").Value = "2026-02-06T12:47:00+01:00"
$ws.Range("J197").Value = ""
$ws.Range("K197").Value = "2026-02-09T13:39:08+00:00"
$ws.Range("L197").Value = "JJD149990121180199217"

# row 198 <- old row 195
$ws.Range("A198").Value = "JJD149990121180199132"
$ws.Range("B198").Value = "dhlparcel-nl"
$ws.Range("C198").Value = "DHL Parcel NL"
$ws.Range("I198").Value = "2026-02-06T12:47:00+01:00"
$ws.Range("J198").Value = ""
$ws.Range("K198").Value = "2026-02-09T13:39:07+00:00"
$ws.Range("L198").Value = "JJD149990121180199132"

# row 199 <- old row 194
$ws.Range("A199").Value = "JJD149990121180198167"
$ws.Range("B199").Value = "dhlparcel-nl"
$ws.Range("C199").Value = "DHL Parcel NL"
$ws.Range("I199").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J199").Value = ""
$ws.Range("K199").Value = "2026-02-09T16:09:17+00:00"
$ws.Range("L199").Value = "JJD149990121180198167"

# row 200 <- old row 201
$ws.Range("A200").Value = "JJD149990121180199072"
$ws.Range("B200").Value = "dhlparcel-nl"
$ws.Range("C200").Value = "DHL Parcel NL"
$ws.Range("I200").Value = "2026-02-06T12:47:00+01:00"
$ws.Range("J200").Value = ""
$ws.Range("K200").Value = "2026-02-09T13:39:10+00:00"
$ws.Range("L200").Value = "JJD149990121180199072"

# row 201 <- old row 196
$ws.Range("A201").Value = "JJD149990121180199090"
$ws.Range("B201").Value = "dhlparcel-nl"
$ws.Range("C201").Value = "DHL Parcel NL"
$ws.Range("I201").Value = "2026-02-09T16:55:00+01:00"
$ws.Range("J201").Value = ""
$ws.Range("K201").Value = "2026-02-09T16:09:31+00:00"
$ws.Range("L201").Value = "JJD149990121180199090"

